# Boolean Peaking and Flexibility Flags.xlsx
#
# Turn off the "peaker" flag for petroleum-fuelled plant types (crude oil and
# heavy/residual fuel oil - HFO) on the "BPaFF-BITPTaP" (Is This Plant Type a
# Peaker) sheet. These cells used to mirror the "natural gas peaker" flag via
# a formula (=B11); they are now hard-set to 0 so HFO/crude oil are no longer
# guaranteed a minimum dispatch level as peakers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPaFF-BITPTaP")

# crude oil (row 15) and heavy or residual fuel oil / HFO (row 16): stop
# tracking the natural gas peaker flag (=B11) and hard-code to "off" (0)
$ws.Range("B15").Value = 0
$ws.Range("B16").Value = 0

# This sheet becomes the active/selected tab, with the given cell selected
$ws.Activate() | Out-Null
$ws.Range("I16").Select() | Out-Null
